$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2916.611
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 3178.5
$ws.Range("K64").Value = 2000
$ws.Range("L64").Value = 3178.5
$ws.Range("M64").Value = -1752
$ws.Range("N64").Value = -3674.5
$ws.Range("H67").Value = 2916.611
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 3178.5
$ws.Range("K67").Value = 2000
$ws.Range("L67").Value = 3178.5
$ws.Range("M67").Value = -1142
$ws.Range("N67").Value = -4894.5
$ws.Range("H74").Value = 4420
$ws.Range("I74").Value = 5360
$ws.Range("J74").Value = 3950
$ws.Range("K74").Value = 5360
$ws.Range("L74").Value = 3950
$ws.Range("M74").Value = -4424
$ws.Range("N74").Value = -5822
$ws.Range("H77").Value = 4420
$ws.Range("I77").Value = 5360
$ws.Range("J77").Value = 3950
$ws.Range("K77").Value = 26800
$ws.Range("L77").Value = 19750
$ws.Range("M77").Value = -22120
$ws.Range("N77").Value = -29110
$ws.Range("H129").Value = 17726.184
$ws.Range("J129").Value = 22956.5
$ws.Range("L129").Value = 68869.5
$ws.Range("N129").Value = -78869.5
$ws.Range("H132").Value = 4466405.5
$ws.Range("I132").Value = 5954862
$ws.Range("J132").Value = 1035.5
$ws.Range("K132").Value = 17864586
$ws.Range("L132").Value = 3106.5
$ws.Range("M132").Value = -17862056
$ws.Range("N132").Value = -8166.5
$ws.Range("H137").Value = 1031.6666
$ws.Range("I137").Value = 801.29034
$ws.Range("J137").Value = 2460
$ws.Range("K137").Value = 2403.87102
$ws.Range("L137").Value = 7380
$ws.Range("M137").Value = 146.12898
$ws.Range("N137").Value = -12480
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 967.39685
$ws.Range("I61").Value = 577.0599999999999
$ws.Range("J61").Value = 2468.6924
$ws.Range("K61").Value = 577.0599999999999
$ws.Range("L61").Value = 2468.6924
$ws.Range("M61").Value = -365.0599999999999
$ws.Range("N61").Value = -2892.6924
$ws.Range("H63").Value = 2001540
$ws.Range("I63").Value = 2501425
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 2501425
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -2500739
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 2001540
$ws.Range("I66").Value = 2501425
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 12507125
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -12503693
$ws.Range("N66").Value = -16864
$ws.Range("H88").Value = 479729.28
$ws.Range("I88").Value = 1432572.2
$ws.Range("J88").Value = 3307.7856
$ws.Range("K88").Value = 1432572.2
$ws.Range("L88").Value = 3307.7856
$ws.Range("M88").Value = -1432166.2
$ws.Range("N88").Value = -4119.7856
$ws.Range("H91").Value = 479729.28
$ws.Range("I91").Value = 1432572.2
$ws.Range("J91").Value = 3307.7856
$ws.Range("K91").Value = 1432572.2
$ws.Range("L91").Value = 3307.7856
$ws.Range("M91").Value = -1431168.2
$ws.Range("N91").Value = -6115.7856
$ws.Range("H136").Value = 967.39685
$ws.Range("I136").Value = 577.0599999999999
$ws.Range("J136").Value = 2468.6924
$ws.Range("K136").Value = 1731.18
$ws.Range("L136").Value = 7406.0772
$ws.Range("M136").Value = 818.8200000000002
$ws.Range("N136").Value = -12506.0772
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2845.5217
$ws.Range("I86").Value = 2269.5
$ws.Range("J86").Value = 3473.9092
$ws.Range("K86").Value = 2269.5
$ws.Range("L86").Value = 3473.9092
$ws.Range("M86").Value = -1146.5
$ws.Range("N86").Value = -5719.9092
$ws.Range("H89").Value = 2845.5217
$ws.Range("I89").Value = 2269.5
$ws.Range("J89").Value = 3473.9092
$ws.Range("K89").Value = 11347.5
$ws.Range("L89").Value = 17369.546
$ws.Range("M89").Value = -5731.5
$ws.Range("N89").Value = -28601.546
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3705942.5
$ws.Range("I31").Value = 1692.931
$ws.Range("J31").Value = 10419895
$ws.Range("K31").Value = 1692.931
$ws.Range("L31").Value = 10419895
$ws.Range("M31").Value = -1397.931
$ws.Range("N31").Value = -10420485
$ws.Range("H34").Value = 3705942.5
$ws.Range("I34").Value = 1692.931
$ws.Range("J34").Value = 10419895
$ws.Range("K34").Value = 1692.931
$ws.Range("L34").Value = 10419895
$ws.Range("M34").Value = -1490.931
$ws.Range("N34").Value = -10420299
$ws.Range("H58").Value = 853.25806
$ws.Range("I58").Value = 765.9545000000001
$ws.Range("J58").Value = 1066.6666
$ws.Range("K58").Value = 765.9545000000001
$ws.Range("L58").Value = 1066.6666
$ws.Range("M58").Value = -562.9545000000001
$ws.Range("N58").Value = -1472.6666
$ws.Range("H132").Value = 3465.3809
$ws.Range("I132").Value = 2565.5833
$ws.Range("J132").Value = 4665.1113
$ws.Range("K132").Value = 7696.749899999999
$ws.Range("L132").Value = 13995.3339
$ws.Range("M132").Value = -5166.749899999999
$ws.Range("N132").Value = -19055.3339
$ws.Range("H134").Value = 971.53125
$ws.Range("I134").Value = 824.04
$ws.Range("J134").Value = 1498.2858
$ws.Range("K134").Value = 2472.12
$ws.Range("L134").Value = 4494.857400000001
$ws.Range("M134").Value = 62.88000000000011
$ws.Range("N134").Value = -9564.857400000001
$ws.Range("H136").Value = 853.25806
$ws.Range("I136").Value = 765.9545000000001
$ws.Range("J136").Value = 1066.6666
$ws.Range("K136").Value = 2297.8635
$ws.Range("L136").Value = 3199.9998
$ws.Range("M136").Value = 252.1364999999996
$ws.Range("N136").Value = -8299.9998
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 10422.223
$ws.Range("I123").Value = 1266.6666
$ws.Range("J123").Value = 15000
$ws.Range("K123").Value = 3799.9998
$ws.Range("L123").Value = 45000
$ws.Range("M123").Value = -1349.9998
$ws.Range("N123").Value = -49900
$ws.Range("H131").Value = 413116.03
$ws.Range("I131").Value = 5582.476
$ws.Range("K131").Value = 16747.428
$ws.Range("M131").Value = -11707.428
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2400
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2400
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 2400
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -4396
$ws.Range("H83").Value = 2400
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2400
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 12000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -21984
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1280
$ws.Range("I68").Value = 1280
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1280
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -531
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 1280
$ws.Range("I71").Value = 1280
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 6400
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -2656
$ws.Range("N71").ClearContents()
$ws.Range("H82").Value = 1931.1111
$ws.Range("I82").Value = 2050
$ws.Range("J82").Value = 980
$ws.Range("K82").Value = 2050
$ws.Range("L82").Value = 980
$ws.Range("M82").Value = -1689
$ws.Range("N82").Value = -1702
$ws.Range("H85").Value = 1931.1111
$ws.Range("I85").Value = 2050
$ws.Range("J85").Value = 980
$ws.Range("K85").Value = 2050
$ws.Range("L85").Value = 980
$ws.Range("M85").Value = -802
$ws.Range("N85").Value = -3476
$ws.Range("H93").Value = 1600.2778
$ws.Range("I93").Value = 1600.5
$ws.Range("J93").Value = 1600
$ws.Range("K93").Value = 1600.5
$ws.Range("L93").Value = 1600
$ws.Range("M93").Value = -352.5
$ws.Range("N93").Value = -4096
$ws.Range("H132").Value = 4898
$ws.Range("I132").Value = 5839.6216
$ws.Range("J132").Value = 2409.4285
$ws.Range("K132").Value = 17518.8648
$ws.Range("L132").Value = 7228.2855
$ws.Range("M132").Value = -14988.8648
$ws.Range("N132").Value = -12288.2855
$ws.Range("H136").Value = 5375.154
$ws.Range("I136").Value = 7430
$ws.Range("J136").Value = 2573.0908
$ws.Range("K136").Value = 22290
$ws.Range("L136").Value = 7719.2724
$ws.Range("M136").Value = -19740
$ws.Range("N136").Value = -12819.2724
$ws.Range("H138").Value = 47532.223
$ws.Range("J138").Value = 47532.223
$ws.Range("L138").Value = 47532.223
$ws.Range("N138").Value = -57812.223
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2875
$ws.Range("I62").Value = 2833.3333
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2833.3333
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -2209.3333
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 2875
$ws.Range("I65").Value = 2833.3333
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 14166.6665
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -11046.6665
$ws.Range("N65").Value = -21240
$ws.Range("H81").Value = 1234.0714
$ws.Range("I81").Value = 1139.75
$ws.Range("J81").Value = 1800
$ws.Range("K81").Value = 2279.5
$ws.Range("L81").Value = 3600
$ws.Range("M81").Value = -1218.5
$ws.Range("N81").Value = -5722
$ws.Range("H84").Value = 1234.0714
$ws.Range("I84").Value = 1139.75
$ws.Range("J84").Value = 1800
$ws.Range("K84").Value = 11397.5
$ws.Range("L84").Value = 18000
$ws.Range("M84").Value = -6093.5
$ws.Range("N84").Value = -28608
